$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Z (column 26), shifting the rest right.
$ws.Columns("Z:Z").Insert()

# Set width of the newly inserted column Z.
$ws.Columns("Z:Z").ColumnWidth = 8.21875

# Populate the new column's header (row 1) and data row (row 2) cells.
$ws.Range("Z1").Value = "STAT"
$ws.Range("Z1").Style = $ws.Range("AA1").Style

$ws.Range("Z2").Value = "K"
$ws.Range("Z2").Style = $ws.Range("AA2").Style

# Update the view state to match the committed selection/scroll position.
$ws.Range("AG8").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("Y1").Column
